$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date strings in column A change format from DD/MM/YYYY to DD-MM-YYYY.
# They must stay literal text (matching the original inlineStr cells), so
# force a text number format before writing, then restore the cell style
# (Excel would otherwise silently reinterpret some of these as real dates).
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    $cell.Style = "Normal"
}

# Updated attendance numbers for rows 3-6 (D=Total, E=Real, F=Duplicate,
# G=Invalid, H=Absent).
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 0
$ws.Cells.Item(3, 7).Value = 1
$ws.Cells.Item(3, 8).Value = 1

$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0

$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0

$ws.Cells.Item(6, 4).Value = 1
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0
